# Insert a new weekly price record as row 22, shifting all existing
# records (previously rows 22-150) down by one row (to rows 23-151).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 22 (pushes rows 22:150 down to 23:151)
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the new record's data
$ws.Range("A22").Value = 1
$ws.Range("B22").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C22").Value = "Arica y Parinacota"
$ws.Range("D22").Value = 44831
$ws.Range("E22").Value = 15
$ws.Range("F22").Value = "Fruta"
$ws.Range("G22").Value = 100108
$ws.Range("H22").Value = "Tropicales y subtropicales"
$ws.Range("I22").Value = 100108002
$ws.Range("J22").Value = "Mango"
$ws.Range("K22").Value = "Sin especificar"
$ws.Range("L22").Value = "Especial"
$ws.Range("M22").Value = 400
$ws.Range("N22").Value = 8000
$ws.Range("O22").Value = 9000
$ws.Range("P22").Value = 8500
$ws.Range("Q22").Value = "`$/bandeja 4 kilos"
$ws.Range("R22").Value = "Brasil"
$ws.Range("S22").Value = 2125
$ws.Range("T22").Value = 4
